$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-7 moves from serial date 45174 (2023-09-05)
# to 45175 (2023-09-06).
$ws.Range("C2:C7").Value = 45175
